# Removing less than USD 5 price from extrapolation calibration because it is just a noise
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3" = 119752.2188536679;  "E3" = -0.01011729461075268; "F3" = 0.2202714586791192;  "G3" = -1.163088336028754;  "H3" = 11.16982911953582
    "D4" = 120604.0864605543;  "E4" = -0.01955736144096514; "F4" = 0.2532773235655477;  "G4" = -1.558699261657695;  "H4" = 14.16771763371706
    "D5" = 121390.320190742;   "E5" = -0.02124510350304856; "F5" = 0.2397519317040472;  "G5" = -0.7043761502860133; "H5" = 7.812699049623616
    "D6" = 121902.7993294967;  "E6" = -0.03264548407403056; "F6" = 0.2718449549411772;  "G6" = -1.149960921263174;  "H6" = 10.10687342162665
    "D8" = 123570.6245224911;  "E8" = -0.04768347842864261; "F8" = 0.2263875357639576;  "G8" = -0.8810990765944922; "H8" = 6.939560630134887
    "D10" = 126384.7604277703; "E10" = -0.1145507063135969; "F10" = 0.441096196145863;  "G10" = -1.938185720727088; "H10" = 9.944337269473078
    "D11" = 128480.9794569601; "E11" = -0.1920957381181317; "F11" = 0.7791410834716266; "G11" = -2.646234149017586; "H11" = 13.02800252663572
    "D12" = 119048.6234017637; "E12" = 0.04062706907592137; "F12" = 0.1562518471311289; "G12" = -1.089849091793474; "H12" = 12.31632432045555
    "D15" = 119011.8400046647; "E15" = 0.08719352372954203; "F15" = 0.1252996970048831; "G15" = -0.1107008630014937; "H15" = 4.691178673692739
    "D17" = 118976.228795658;  "E17" = 0.0536154588118516;  "F17" = 0.1645618851359608; "G17" = -1.682696326551386;  "H17" = 15.80779403450412
    "D19" = 119054.8228035294; "E19" = 0.01508357598063807; "F19" = 0.1598502418917181; "G19" = -1.050369021455628;  "H19" = 10.31880630607952
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
